# Update cryptos list: refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.513.90"
$ws.Range("E2").Value = "  -0.10%  "

$ws.Range("D3").Value = "1.901.78"
$ws.Range("E3").Value = "  +1.35%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.04"
$ws.Range("E5").Value = "  +1.21%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4915"
$ws.Range("E7").Value = "  +0.81%  "

$ws.Range("E8").Value = "  +1.30%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06701"
$ws.Range("E9").Value = "  +0.52%  "

$ws.Range("D10").Value = "1.898.20"
$ws.Range("E10").Value = "  +1.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.09"
$ws.Range("E11").Value = "  +3.13%  "

$ws.Range("E12").Value = "  +1.34%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.177"
$ws.Range("E13").Value = "  +3.65%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "88.07"
$ws.Range("E14").Value = "  -0.61%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6693"
$ws.Range("E15").Value = "  +3.01%  "

$ws.Range("D16").Value = "30.492.27"
$ws.Range("E16").Value = "  +0.01%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.50"
$ws.Range("E17").Value = "  +3.99%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007896"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9994"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.450"
$ws.Range("E20").Value = "  +15.63%  "

$ws.Range("D21").Value = "2.143.67"
$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "196.35"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.135"
$ws.Range("E24").Value = "  +0.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.525"
$ws.Range("E25").Value = "  +1.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.76"
$ws.Range("E26").Value = "  +3.61%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.39"
$ws.Range("E27").Value = "  -0.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.945"
$ws.Range("E28").Value = "  +6.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.482"
$ws.Range("E29").Value = "  +5.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.355"
$ws.Range("E30").Value = "  +2.51%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09180"
$ws.Range("E31").Value = "  +1.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.097"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05173"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7452"
$ws.Range("E34").Value = "  +3.39%  "

$ws.Range("E35").Value = "  +2.90%  "

$ws.Range("E36").Value = "  +1.03%  "

$ws.Range("E37").Value = "  +0.31%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.684"
$ws.Range("E38").Value = "  +0.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9268"
$ws.Range("E39").Value = "  +0.74%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.068"
$ws.Range("E40").Value = "  +1.30%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4406"
$ws.Range("E41").Value = "  +0.52%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.944"
$ws.Range("E42").Value = "  +4.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.01"
$ws.Range("E43").Value = "  +2.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.55"
$ws.Range("E44").Value = "  +22.11%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9960"
$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1374"
$ws.Range("E46").Value = "  +3.49%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.612"
$ws.Range("E47").Value = "  +3.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.026"
$ws.Range("E48").Value = "  +4.82%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.11"
$ws.Range("E49").Value = "  +6.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05842"
$ws.Range("E50").Value = "  +0.30%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3938"
$ws.Range("E51").Value = "  -1.97%  "
